# Insert a new row above the existing "Width:" row (row 8) in the
# "Survey Table Properties" block, adding a pseudo-random question width
# option. Everything below shifts down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 8 ("Width:" ...) and everything after it down by one row.
$ws.Rows.Item(8).Insert()

# Populate the new row 8, matching the look of its sibling label cells
# (column A holds the label text, column B is the blank input cell).
$ws.Range("A8").Value = "Pseudo-Random Question Width:"
$ws.Range("A8").Style = "Normal"

# Leave the selection on B8, matching where the user clicked after typing
# the new label.
$ws.Range("B8").Select()
